# Update NATMI TPM output with refreshed numbers.
# - Rows 2-10 (FAPs/MuSCs/Neutrophils x ECs/FAPs/MuSCs) get refreshed metrics.
# - Rows 11-13 (target cluster "Inflammatory-Mac") are removed entirely, since
#   that cluster no longer appears in the refreshed TPM run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..10 (columns A..T)
$data = @(
    @("FAPs","Artn","Gfra1","ECs",2,0.6666666666666666,0.6226999999999999,1.8681,0.8812810914468889,0.881281091446889,2,0.6666666666666666,0.1389376666666667,0.416813,0.01722256533596611,0.01722256533596611,0.08651648503333331,0.7786483652999999,0.01517792117679557,0.01517792117679557),
    @("FAPs","Artn","Gfra1","FAPs",2,0.6666666666666666,0.6226999999999999,1.8681,0.8812810914468889,0.881281091446889,3,1,4.220039333333333,12.660118,0.5231115858095611,0.5231115858095611,2.627818492866666,23.6503664358,0.4610083492907629,0.461008349290763),
    @("FAPs","Artn","Gfra1","MuSCs",2,0.6666666666666666,0.6226999999999999,1.8681,0.8812810914468889,0.881281091446889,3,1,3.708210666666667,11.124632,0.4596658488544727,0.4596658488544728,2.309102782133333,20.7819250392,0.4050948209793304,0.4050948209793305),
    @("MuSCs","Artn","Gfra1","ECs",1,0.3333333333333333,0.04989433333333334,0.149683,0.07061334918422178,0.07061334918422178,2,0.6666666666666666,0.1389376666666667,0.416813,0.01722256533596611,0.01722256533596611,0.006932202253222222,0.062389820279,0.001216143019916649,0.001216143019916649),
    @("MuSCs","Artn","Gfra1","FAPs",1,0.3333333333333333,0.04989433333333334,0.149683,0.07061334918422178,0.07061334918422178,3,1,4.220039333333333,12.660118,0.5231115858095611,0.5231115858095611,0.2105560491771111,1.895004442594,0.03693866107108253,0.03693866107108253),
    @("MuSCs","Artn","Gfra1","MuSCs",1,0.3333333333333333,0.04989433333333334,0.149683,0.07061334918422178,0.07061334918422178,3,1,3.708210666666667,11.124632,0.4596658488544727,0.4596658488544728,0.1850186990728889,1.665168291656,0.03245854509322259,0.03245854509322259),
    @("Neutrophils","Artn","Gfra1","ECs",1,0.3333333333333333,0.03399066666666666,0.101972,0.04810555936888933,0.04810555936888933,2,0.6666666666666666,0.1389376666666667,0.416813,0.01722256533596611,0.01722256533596611,0.00472258391511111,0.04250325523599999,0.000828501139253893,0.0008285011392538933),
    @("Neutrophils","Artn","Gfra1","FAPs",1,0.3333333333333333,0.03399066666666666,0.101972,0.04810555936888933,0.04810555936888933,3,1,4.220039333333333,12.660118,0.5231115858095611,0.5231115858095611,0.1434419502995555,1.290977552696,0.02516457544771569,0.02516457544771569),
    @("Neutrophils","Artn","Gfra1","MuSCs",1,0.3333333333333333,0.03399066666666666,0.101972,0.04810555936888933,0.04810555936888933,3,1,3.708210666666667,11.124632,0.4596658488544727,0.4596658488544728,0.1260445527004444,1.134400974304,0.02211248278191974,0.02211248278191975)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = 1 + $c
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}

# Remove the now-obsolete rows that referenced the "Inflammatory-Mac" target
# cluster (old rows 11, 12, 13) so the sheet shrinks from 13 to 10 data rows.
$ws.Rows.Item(11).Resize(3).EntireRow.Delete()
